$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 8.84 = 36862.48 pesos`n✅ 36862.48 pesos = 8.84 = 957.28 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 113.15
$wsTasas.Range("O10").Value = 4170.99
$wsTasas.Range("N12").Value = 4170
$wsTasas.Range("O12").Value = 108.29
